# Applies the automatic-update edit: rows 10-21 in the sheet get their
# observation records reshuffled among the row slots (rows 16 and 20 keep
# their original content as anchors). Only the affected cells are rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (was id 111936795 before shuffle source)
$ws.Range("A10").Value = 111936795
$ws.Range("B10").Value = 56398
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 448749.3706757246
$ws.Range("R10").Value = 7087421.839990681

# Row 11 (was id 111936796 before shuffle source)
$ws.Range("A11").Value = 111936796
$ws.Range("B11").Value = 56398
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 448882.8980770012
$ws.Range("R11").Value = 7087229.443335658

# Row 12 (was id 111936868 before shuffle source)
$ws.Range("A12").Value = 111936868
$ws.Range("B12").Value = 89423
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 448988.017639213
$ws.Range("R12").Value = 7087186.778340456

# Row 13 (was id 111936867 before shuffle source)
$ws.Range("A13").Value = 111936867
$ws.Range("B13").Value = 89423
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = "Granticka"
$ws.Range("G13").Value = "Porodaedalea chrysoloma"
$ws.Range("H13").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 448791.554596175
$ws.Range("R13").Value = 7087386.366048628

# Row 14 (was id 111936870 before shuffle source)
$ws.Range("A14").Value = 111936870
$ws.Range("B14").Value = 89423
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 5432
$ws.Range("F14").Value = "Granticka"
$ws.Range("G14").Value = "Porodaedalea chrysoloma"
$ws.Range("H14").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q14").Value = 449019.027096529
$ws.Range("R14").Value = 7087276.979166135

# Row 15 (was id 111936866 before shuffle source)
$ws.Range("A15").Value = 111936866
$ws.Range("B15").Value = 89423
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 5432
$ws.Range("F15").Value = "Granticka"
$ws.Range("G15").Value = "Porodaedalea chrysoloma"
$ws.Range("H15").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q15").Value = 448765.5992023234
$ws.Range("R15").Value = 7087416.731054713

# Row 17 (was id 111936869 before shuffle source)
$ws.Range("A17").Value = 111936869
$ws.Range("B17").Value = 89423
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 5432
$ws.Range("F17").Value = "Granticka"
$ws.Range("G17").Value = "Porodaedalea chrysoloma"
$ws.Range("H17").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q17").Value = 449143.8568242944
$ws.Range("R17").Value = 7087117.752608996

# Row 18 (was id 111936893 before shuffle source)
$ws.Range("A18").Value = 111936893
$ws.Range("B18").Value = 77515
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("Q18").Value = 448742.3011697636
$ws.Range("R18").Value = 7087501.648173723

# Row 19 (was id 111936858 before shuffle source)
$ws.Range("A19").Value = 111936858
$ws.Range("B19").Value = 89845
$ws.Range("D19").Value = "VU"
$ws.Range("E19").Value = 1209
$ws.Range("F19").Value = "Rynkskinn"
$ws.Range("G19").Value = "Phlebia centrifuga"
$ws.Range("H19").Value = "P.Karst."
$ws.Range("Q19").Value = 448737.3665225056
$ws.Range("R19").Value = 7087496.445579056

# Row 21 (was id 111936792 before shuffle source)
$ws.Range("A21").Value = 111936792
$ws.Range("B21").Value = 90087
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 3298
$ws.Range("F21").Value = "Trådticka"
$ws.Range("G21").Value = "Climacocystis borealis"
$ws.Range("H21").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q21").Value = 448761.1990147882
$ws.Range("R21").Value = 7087578.827763715

# Row 11 becomes a "Tretåig hackspett" bird record: it needs the bird-only
# Ålder-Stadium/Kön/Aktivitet/Metod placeholder cells (left blank, same as
# the other bird rows) and the public comment.
$ws.Range("K11").Value = "'"
$ws.Range("L11").Value = "'"
$ws.Range("M11").Value = "'"
$ws.Range("N11").Value = "'"
$ws.Range("AC11").Value = "ringhack äldre"

# Row 14 becomes a "Granticka" fungus record: it no longer needs the
# bird-only placeholder cells or the public comment, so clear them out.
$ws.Range("K14:N14").ClearContents()
$ws.Range("AC14").ClearContents()
